# The commit swaps the deck's applied colour theme from the "Integral"
# palette over to the stock "Office Theme" palette (the colours that used
# to live only on the orphaned notes-master theme part). Re-apply every
# theme colour slot (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) on the
# slide master's theme so the whole deck renders with the Office Theme
# colours instead of the Integral ones.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# RGB() long values for the Office Theme colour scheme, in the standard
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order.
$officeThemeColors = @(
    0,        # dk1      000000
    16777215, # lt1      FFFFFF
    6968388,  # dk2      44546A
    15132391, # lt2      E7E6E6
    13998939, # accent1  5B9BD5
    3243501,  # accent2  ED7D31
    10855845, # accent3  A5A5A5
    49407,    # accent4  FFC000
    12874308, # accent5  4472C4
    4697456,  # accent6  70AD47
    12673797, # hlink    0563C1
    7491477   # folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
